# Apply targeted cell value updates to Sheet1 (odds data refresh)
# Generated from the authoritative cell-level diff between before/after states.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: columns J, K
$ws.Cells.Item(2, 10).Value = 1.03   # J2
$ws.Cells.Item(2, 11).Value = 15   # K2

# Row 3: columns N, O
$ws.Cells.Item(3, 14).Value = 2.15   # N3
$ws.Cells.Item(3, 15).Value = 1.67   # O3

# Row 11: columns G, H, I, J, K, W, X, AA, AH
$ws.Cells.Item(11, 7).Value = 5.75   # G11
$ws.Cells.Item(11, 8).Value = 3.7   # H11
$ws.Cells.Item(11, 9).Value = 1.57   # I11
$ws.Cells.Item(11, 10).Value = 1.07   # J11
$ws.Cells.Item(11, 11).Value = 9   # K11
$ws.Cells.Item(11, 23).Value = 67   # W11
$ws.Cells.Item(11, 24).Value = 51   # X11
$ws.Cells.Item(11, 27).Value = 7.5   # AA11
$ws.Cells.Item(11, 34).Value = 11   # AH11

# Row 12: columns G, H, I, N, O, T, U, X, Z, AD, AF, AG, AI
$ws.Cells.Item(12, 7).Value = 1.73   # G12
$ws.Cells.Item(12, 8).Value = 3.5   # H12
$ws.Cells.Item(12, 9).Value = 4   # I12
$ws.Cells.Item(12, 14).Value = 1.7   # N12
$ws.Cells.Item(12, 15).Value = 2.1   # O12
$ws.Cells.Item(12, 20).Value = 8.5   # T12
$ws.Cells.Item(12, 21).Value = 9   # U12
$ws.Cells.Item(12, 24).Value = 13   # X12
$ws.Cells.Item(12, 26).Value = 13   # Z12
$ws.Cells.Item(12, 30).Value = 151   # AD12
$ws.Cells.Item(12, 32).Value = 23   # AF12
$ws.Cells.Item(12, 33).Value = 15   # AG12
$ws.Cells.Item(12, 35).Value = 34   # AI12

# Row 13: columns G, H, I, L, M, N, O, U, V, W, X, Y, Z, AE, AF, AH
$ws.Cells.Item(13, 7).Value = 3.3   # G13
$ws.Cells.Item(13, 8).Value = 3.4   # H13
$ws.Cells.Item(13, 9).Value = 1.9   # I13
$ws.Cells.Item(13, 12).Value = 1.18   # L13
$ws.Cells.Item(13, 13).Value = 4.5   # M13
$ws.Cells.Item(13, 14).Value = 1.65   # N13
$ws.Cells.Item(13, 15).Value = 2.2   # O13
$ws.Cells.Item(13, 21).Value = 21   # U13
$ws.Cells.Item(13, 22).Value = 13   # V13
$ws.Cells.Item(13, 23).Value = 41   # W13
$ws.Cells.Item(13, 24).Value = 26   # X13
$ws.Cells.Item(13, 25).Value = 29   # Y13
$ws.Cells.Item(13, 26).Value = 15   # Z13
$ws.Cells.Item(13, 31).Value = 9.5   # AE13
$ws.Cells.Item(13, 32).Value = 11   # AF13
$ws.Cells.Item(13, 34).Value = 17   # AH13

# Row 14: columns G, H, I, J, K, N, O, Z, AE
$ws.Cells.Item(14, 7).Value = 1.53   # G14
$ws.Cells.Item(14, 8).Value = 3.6   # H14
$ws.Cells.Item(14, 9).Value = 5.25   # I14
$ws.Cells.Item(14, 10).Value = 1.04   # J14
$ws.Cells.Item(14, 11).Value = 13   # K14
$ws.Cells.Item(14, 14).Value = 1.75   # N14
$ws.Cells.Item(14, 15).Value = 2.05   # O14
$ws.Cells.Item(14, 26).Value = 12   # Z14
$ws.Cells.Item(14, 31).Value = 17   # AE14

# Row 15: columns G, I, U, W, X, AG
$ws.Cells.Item(15, 7).Value = 1.95   # G15
$ws.Cells.Item(15, 9).Value = 3.3   # I15
$ws.Cells.Item(15, 21).Value = 10   # U15
$ws.Cells.Item(15, 23).Value = 19   # W15
$ws.Cells.Item(15, 24).Value = 17   # X15
$ws.Cells.Item(15, 33).Value = 12   # AG15

# Row 21: columns G, H, I, L, M, N, O, P, Q, R, S, T, U, V, W, X, Y, Z, AA, AB, AC, AE, AF, AG, AH, AI, AJ
$ws.Cells.Item(21, 7).Value = 1.47   # G21
$ws.Cells.Item(21, 8).Value = 3.65   # H21
$ws.Cells.Item(21, 9).Value = 7.5   # I21
$ws.Cells.Item(21, 12).Value = 1.44   # L21
$ws.Cells.Item(21, 13).Value = 2.4   # M21
$ws.Cells.Item(21, 14).Value = 2.27   # N21
$ws.Cells.Item(21, 15).Value = 1.5   # O21
$ws.Cells.Item(21, 16).Value = 1.52   # P21
$ws.Cells.Item(21, 17).Value = 2.2   # Q21
$ws.Cells.Item(21, 18).Value = 2.45   # R21
$ws.Cells.Item(21, 19).Value = 1.42   # S21
$ws.Cells.Item(21, 20).Value = 4.6   # T21
$ws.Cells.Item(21, 21).Value = 5.4   # U21
$ws.Cells.Item(21, 22).Value = 9.25   # V21
$ws.Cells.Item(21, 23).Value = 9.25   # W21
$ws.Cells.Item(21, 24).Value = 15.5   # X21
$ws.Cells.Item(21, 25).Value = 50   # Y21
$ws.Cells.Item(21, 26).Value = 6.8   # Z21
$ws.Cells.Item(21, 27).Value = 7.8   # AA21
$ws.Cells.Item(21, 28).Value = 30   # AB21
$ws.Cells.Item(21, 29).Value = 250   # AC21
$ws.Cells.Item(21, 31).Value = 13   # AE21
$ws.Cells.Item(21, 32).Value = 45   # AF21
$ws.Cells.Item(21, 33).Value = 28   # AG21
$ws.Cells.Item(21, 34).Value = 250   # AH21
$ws.Cells.Item(21, 35).Value = 150   # AI21
$ws.Cells.Item(21, 36).Value = 150   # AJ21

# Row 25: columns G, I, P, Q, U, V, W, AD, AE, AF, AG, AH, AI, AJ
$ws.Cells.Item(25, 7).Value = 1.78   # G25
$ws.Cells.Item(25, 9).Value = 4.25   # I25
$ws.Cells.Item(25, 16).Value = 1.39   # P25
$ws.Cells.Item(25, 17).Value = 2.42   # Q25
$ws.Cells.Item(25, 21).Value = 6.6   # U25
$ws.Cells.Item(25, 22).Value = 7.1   # V25
$ws.Cells.Item(25, 23).Value = 11.75   # W25
$ws.Cells.Item(25, 30).Value = 500   # AD25
$ws.Cells.Item(25, 31).Value = 8.5   # AE25
$ws.Cells.Item(25, 32).Value = 18   # AF25
$ws.Cells.Item(25, 33).Value = 12   # AG25
$ws.Cells.Item(25, 34).Value = 55   # AH25
$ws.Cells.Item(25, 35).Value = 35   # AI25
$ws.Cells.Item(25, 36).Value = 40   # AJ25

# Row 27: columns N, O
$ws.Cells.Item(27, 14).Value = 1.73   # N27
$ws.Cells.Item(27, 15).Value = 2.08   # O27

# Row 28: columns N, T, U, AE, AF, AJ
$ws.Cells.Item(28, 14).Value = 2   # N28
$ws.Cells.Item(28, 20).Value = 8   # T28
$ws.Cells.Item(28, 21).Value = 13   # U28
$ws.Cells.Item(28, 31).Value = 8.5   # AE28
$ws.Cells.Item(28, 32).Value = 14.5   # AF28
$ws.Cells.Item(28, 36).Value = 30   # AJ28

# Row 29: columns G, I, K, R, T, U, V, W, X, Z, AB, AE, AF, AG, AH, AI
$ws.Cells.Item(29, 7).Value = 2.15   # G29
$ws.Cells.Item(29, 9).Value = 3.95   # I29
$ws.Cells.Item(29, 11).Value = 5.6   # K29
$ws.Cells.Item(29, 18).Value = 1.88   # R29
$ws.Cells.Item(29, 20).Value = 6.1   # T29
$ws.Cells.Item(29, 21).Value = 9.75   # U29
$ws.Cells.Item(29, 22).Value = 8.75   # V29
$ws.Cells.Item(29, 23).Value = 22   # W29
$ws.Cells.Item(29, 24).Value = 19.5   # X29
$ws.Cells.Item(29, 26).Value = 5.6   # Z29
$ws.Cells.Item(29, 28).Value = 14   # AB29
$ws.Cells.Item(29, 31).Value = 9.5   # AE29
$ws.Cells.Item(29, 32).Value = 22   # AF29
$ws.Cells.Item(29, 33).Value = 13   # AG29
$ws.Cells.Item(29, 34).Value = 70   # AH29
$ws.Cells.Item(29, 35).Value = 40   # AI29

# Row 31: columns G, I, T, U, W, AC, AE, AF, AG, AH, AJ
$ws.Cells.Item(31, 7).Value = 1.8   # G31
$ws.Cells.Item(31, 9).Value = 3.9   # I31
$ws.Cells.Item(31, 20).Value = 7.8   # T31
$ws.Cells.Item(31, 21).Value = 9   # U31
$ws.Cells.Item(31, 23).Value = 14.5   # W31
$ws.Cells.Item(31, 29).Value = 60   # AC31
$ws.Cells.Item(31, 31).Value = 12   # AE31
$ws.Cells.Item(31, 32).Value = 22   # AF31
$ws.Cells.Item(31, 33).Value = 13   # AG31
$ws.Cells.Item(31, 34).Value = 55   # AH31
$ws.Cells.Item(31, 36).Value = 37   # AJ31

# Row 32: columns G, H, I, K, O, Q, R, S, T, U, V, W, X, Y, Z, AA, AB, AC, AD, AE, AF, AG, AH, AI, AJ
$ws.Cells.Item(32, 7).Value = 1.65   # G32
$ws.Cells.Item(32, 8).Value = 3.6   # H32
$ws.Cells.Item(32, 9).Value = 5   # I32
$ws.Cells.Item(32, 11).Value = 7.3   # K32
$ws.Cells.Item(32, 15).Value = 1.85   # O32
$ws.Cells.Item(32, 17).Value = 2.67   # Q32
$ws.Cells.Item(32, 18).Value = 1.83   # R32
$ws.Cells.Item(32, 19).Value = 1.88   # S32
$ws.Cells.Item(32, 20).Value = 6.8   # T32
$ws.Cells.Item(32, 21).Value = 7.8   # U32
$ws.Cells.Item(32, 22).Value = 8   # V32
$ws.Cells.Item(32, 23).Value = 13   # W32
$ws.Cells.Item(32, 24).Value = 13   # X32
$ws.Cells.Item(32, 25).Value = 25   # Y32
$ws.Cells.Item(32, 26).Value = 7.3   # Z32
$ws.Cells.Item(32, 27).Value = 7   # AA32
$ws.Cells.Item(32, 28).Value = 16   # AB32
$ws.Cells.Item(32, 29).Value = 75   # AC32
$ws.Cells.Item(32, 30).Value = 600   # AD32
$ws.Cells.Item(32, 31).Value = 12.5   # AE32
$ws.Cells.Item(32, 32).Value = 29   # AF32
$ws.Cells.Item(32, 33).Value = 16   # AG32
$ws.Cells.Item(32, 34).Value = 90   # AH32
$ws.Cells.Item(32, 35).Value = 55   # AI32
$ws.Cells.Item(32, 36).Value = 55   # AJ32

# Row 33: columns G, I, S, T, U, V, W, X, AE, AF, AG, AH, AI
$ws.Cells.Item(33, 7).Value = 2.47   # G33
$ws.Cells.Item(33, 9).Value = 3.1   # I33
$ws.Cells.Item(33, 19).Value = 1.7   # S33
$ws.Cells.Item(33, 20).Value = 6.4   # T33
$ws.Cells.Item(33, 21).Value = 11   # U33
$ws.Cells.Item(33, 22).Value = 9.75   # V33
$ws.Cells.Item(33, 23).Value = 27   # W33
$ws.Cells.Item(33, 24).Value = 24   # X33
$ws.Cells.Item(33, 31).Value = 7.1   # AE33
$ws.Cells.Item(33, 32).Value = 14.5   # AF33
$ws.Cells.Item(33, 33).Value = 11.5   # AG33
$ws.Cells.Item(33, 34).Value = 40   # AH33
$ws.Cells.Item(33, 35).Value = 35   # AI33

# Row 34: columns K, R, S, AB, AD, AG, AI
$ws.Cells.Item(34, 11).Value = 7.5   # K34
$ws.Cells.Item(34, 18).Value = 2   # R34
$ws.Cells.Item(34, 19).Value = 1.75   # S34
$ws.Cells.Item(34, 28).Value = 17   # AB34
$ws.Cells.Item(34, 30).Value = 451   # AD34
$ws.Cells.Item(34, 33).Value = 13   # AG34
$ws.Cells.Item(34, 35).Value = 34   # AI34

# Row 36: columns G, I, U, V, Y, AE, AJ
$ws.Cells.Item(36, 7).Value = 2.38   # G36
$ws.Cells.Item(36, 9).Value = 3   # I36
$ws.Cells.Item(36, 21).Value = 13   # U36
$ws.Cells.Item(36, 22).Value = 9.5   # V36
$ws.Cells.Item(36, 25).Value = 21   # Y36
$ws.Cells.Item(36, 31).Value = 13   # AE36
$ws.Cells.Item(36, 36).Value = 26   # AJ36

# Row 37: columns G, H, I, N, O, R, S, T, U, V, W, X, Y, Z, AA, AB, AC, AE, AF, AG, AI, AJ
$ws.Cells.Item(37, 7).Value = 1.08   # G37
$ws.Cells.Item(37, 8).Value = 7.7   # H37
$ws.Cells.Item(37, 9).Value = 29   # I37
$ws.Cells.Item(37, 14).Value = 1.34   # N37
$ws.Cells.Item(37, 15).Value = 2.77   # O37
$ws.Cells.Item(37, 18).Value = 2.47   # R37
$ws.Cells.Item(37, 19).Value = 1.42   # S37
$ws.Cells.Item(37, 20).Value = 9.25   # T37
$ws.Cells.Item(37, 21).Value = 6   # U37
$ws.Cells.Item(37, 22).Value = 13   # V37
$ws.Cells.Item(37, 23).Value = 5.7   # W37
$ws.Cells.Item(37, 24).Value = 12   # X37
$ws.Cells.Item(37, 25).Value = 50   # Y37
$ws.Cells.Item(37, 26).Value = 18.5   # Z37
$ws.Cells.Item(37, 27).Value = 21   # AA37
$ws.Cells.Item(37, 28).Value = 50   # AB37
$ws.Cells.Item(37, 29).Value = 250   # AC37
$ws.Cells.Item(37, 31).Value = 100   # AE37
$ws.Cells.Item(37, 32).Value = 500   # AF37
$ws.Cells.Item(37, 33).Value = 120   # AG37
$ws.Cells.Item(37, 35).Value = 900   # AI37
$ws.Cells.Item(37, 36).Value = 400   # AJ37

# Row 38: columns H, I
$ws.Cells.Item(38, 8).Value = 6.5   # H38
$ws.Cells.Item(38, 9).Value = 1.2   # I38

# Row 40: columns G, I, J, K, N, O, R, S, Z, AB, AI
$ws.Cells.Item(40, 7).Value = 1.9   # G40
$ws.Cells.Item(40, 9).Value = 3.6   # I40
$ws.Cells.Item(40, 10).Value = 23   # J40
$ws.Cells.Item(40, 11).Value = 1.03   # K40
$ws.Cells.Item(40, 14).Value = 1.36   # N40
$ws.Cells.Item(40, 15).Value = 3   # O40
$ws.Cells.Item(40, 18).Value = 1.36   # R40
$ws.Cells.Item(40, 19).Value = 3   # S40
$ws.Cells.Item(40, 26).Value = 23   # Z40
$ws.Cells.Item(40, 28).Value = 10   # AB40
$ws.Cells.Item(40, 35).Value = 23   # AI40

# Row 46: columns G, K, P, Q, S, T, U, W, X, Y, Z, AA, AE, AF, AI, AJ
$ws.Cells.Item(46, 7).Value = 2.5   # G46
$ws.Cells.Item(46, 11).Value = 7.4   # K46
$ws.Cells.Item(46, 16).Value = 1.4   # P46
$ws.Cells.Item(46, 17).Value = 2.72   # Q46
$ws.Cells.Item(46, 19).Value = 2.18   # S46
$ws.Cells.Item(46, 20).Value = 9.25   # T46
$ws.Cells.Item(46, 21).Value = 13.5   # U46
$ws.Cells.Item(46, 23).Value = 28   # W46
$ws.Cells.Item(46, 24).Value = 19.5   # X46
$ws.Cells.Item(46, 25).Value = 25   # Y46
$ws.Cells.Item(46, 26).Value = 7.4   # Z46
$ws.Cells.Item(46, 27).Value = 6.2   # AA46
$ws.Cells.Item(46, 31).Value = 9.75   # AE46
$ws.Cells.Item(46, 32).Value = 15   # AF46
$ws.Cells.Item(46, 35).Value = 22   # AI46
$ws.Cells.Item(46, 36).Value = 27   # AJ46

Write-Host "Applied 236 cell updates"
